# Apply the "break out stock.yaml completed" update:
#  - Q56: 2 -> 0
#  - O457: 0 -> 2
#  - R459, R460: blank -> 0
#  - Append new weekly rows 461-473 (A:Q), R left blank (not yet backfilled)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- point fixes on existing rows ---
$ws.Cells.Item(56, 17).Value = 0     # Q56
$ws.Cells.Item(457, 15).Value = 2    # O457
$ws.Cells.Item(459, 18).Value = 0    # R459
$ws.Cells.Item(460, 18).Value = 0    # R460

# --- new rows 461-473, columns A..Q (1..17) ---
$newRows = @(
    @(461, 45474, 628.4500122070312, 663.4000244140625, 627.4500122070312, 644.4000244140625, 638.9318237304688,  8130533, 2024, 7,  1, 0, 0, 0, 27, 0, 0, 0),
    @(462, 45481, 649,               662.7999877929688, 633.5999755859375, 634.7000122070312, 629.3140869140625,  4045867, 2024, 7,  8, 0, 0, 0, 28, 0, 0, 0),
    @(463, 45488, 634.7000122070312, 651.5499877929688, 612.4000244140625, 619.2000122070312, 613.9456176757812,  5008436, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 1),
    @(464, 45495, 612,               656.9000244140625, 612,               648.4500122070312, 642.9474487304688,  5317205, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(465, 45502, 649.2000122070312, 687.2999877929688, 646.5499877929688, 656.9000244140625, 651.3257446289062,  7358752, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
    @(466, 45509, 648.5,             665,               598.0999755859375, 610.7000122070312, 605.5177612304688,  6187201, 2024, 8,  5, 0, 0, 0, 32, 0, 0, 0),
    @(467, 45516, 605,               624.9500122070312, 584.3499755859375, 601.2000122070312, 596.098388671875,   3191630, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0),
    @(468, 45523, 602,               605.8499755859375, 584,               594,               588.95947265625,    6241230, 2024, 8, 19, 0, 0, 0, 34, 2, 0, 0),
    @(469, 45530, 597.7000122070312, 615.9000244140625, 588.5999755859375, 606.9500122070312, 601.7996215820312,  7451055, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0),
    @(470, 45537, 637.1500244140625, 689.9500122070312, 630.75,            667,               661.3400268554688, 61919766, 2024, 9,  2, 0, 0, 0, 36, 0, 0, 0),
    @(471, 45544, 664.7999877929688, 665.4500122070312, 632,               633,               627.6285400390625,  5749293, 2024, 9,  9, 0, 0, 0, 37, 0, 0, 0),
    @(472, 45551, 632.9500122070312, 641.4500122070312, 608.9000244140625, 614.2000122070312, 614.2000122070312,  4285357, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0),
    @(473, 45558, 614.9500122070312, 623.7000122070312, 597,               610.7999877929688, 610.7999877929688,  5518075, 2024, 9, 23, 0, 0, 0, 39, 0, 0, 0)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    for ($col = 1; $col -le 17; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $r[$col]
    }
    # Datetime column keeps the same date/time number format as the rest of column A
    $ws.Cells.Item($rowNum, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Host "Edit complete."
